# Add a new "2021" column (column Y) to the worksheet, mirroring the
# formatting of the existing "2020" column (column X), and populate the
# new 2021 data values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for the new column Y (2021), keyed by row number.
$values = @{
    4  = 2021
    5  = 46.69
    6  = 52.52
    7  = 43.22
    8  = 51.31
    9  = 41.31
    10 = 52.43
    11 = 49.27
    12 = 31.68
    13 = 35.590000000000003
    14 = 55.28
    15 = 61.02
    16 = 48.72
}

foreach ($row in $values.Keys) {
    $srcCell = $ws.Range("X$row")
    $dstCell = $ws.Range("Y$row")

    # Write the value first.
    $dstCell.Value = $values[$row]

    # Mirror the source (2020) cell's formatting onto the new (2021) cell.
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
    $dstCell.Value = $values[$row]
}

# Update the sheet view so the new column is reflected in the selection /
# scroll position, matching the author's final view state.
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("AA15").Select()
